$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 14): month abbreviations + sum
$months = @("jan","feb","mar","apr","may","jun","jul","aug","sep","oct","nov","dec")
for ($i = 0; $i -lt $months.Length; $i++) {
    $col = 3 + $i   # C = 3
    $ws.Cells.Item(14, $col).Value = $months[$i]
}
$ws.Cells.Item(14, 15).Value = "sum"   # O14 = "sum"

# Data row (row 15): days per month, left-aligned
$days = @(31,28,31,30,31,30,31,31,30,31,30,31)
for ($i = 0; $i -lt $days.Length; $i++) {
    $col = 3 + $i
    $cell = $ws.Cells.Item(15, $col)
    $cell.Value = $days[$i]
    $cell.HorizontalAlignment = -4131   # xlHAlignLeft
}

# O15 = sum formula
$ws.Cells.Item(15, 15).Formula = "=SUM(C15:N15)"

# Update selection to O16
$ws.Range("O16").Select()
